# fixed some bugs in scaleDownEnding
# The test data rows (symbol, reel1..reel5) on the worksheet get reordered
# to reflect the corrected scaleDownEnding behavior. Row 1 (headers) and
# rows 22-26 (per-reel totals) are unchanged; only rows 2-21 (columns A:F)
# are rewritten in their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(301, 6, 45, 30, 60, 45),
    @(801, 3, 67, 65, 52, 45),
    @(1203, 3, 15, 15, 15, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(701, 3, 90, 45, 97, 15),
    @(101, 9, 30, 15, 60, 15),
    @(902, 1, 0, 0, 0, 0),
    @(501, 9, 52, 30, 75, 45),
    @(401, 9, 48, 67, 75, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(901, 16, 15, 45, 60, 60),
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(1201, 2, 10, 10, 10, 10),
    @(1, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(502, 0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $row = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}
